$wb = $excel.ActiveWorkbook

# --- Update the Version value on the isa_template sheet ---
$templateSheet = $wb.Worksheets.Item("isa_template")
$templateSheet.Range("B4").Value = "1.0.2"

# --- Update the MIAPPE ontology related headers on the GasExchange sheet/table ---
$dataSheet = $wb.Worksheets.Item("GasExchange")
$dataSheet.Range("M1").Value = "Characteristic [Sample information]"
$dataSheet.Range("N1").Value = "Term Source REF (MIAPPE:0178)"
$dataSheet.Range("O1").Value = "Term Accession Number (MIAPPE:0178)"

# Keep the underlying ListObject (table) column names in sync as well,
# in case Excel does not auto-propagate the header-row edits above.
$table = $dataSheet.ListObjects.Item(1)
$table.ListColumns.Item(13).Name = "Characteristic [Sample information]"
$table.ListColumns.Item(14).Name = "Term Source REF (MIAPPE:0178)"
$table.ListColumns.Item(15).Name = "Term Accession Number (MIAPPE:0178)"
